# Scheduled-runner style update: refresh computed market-price / profit
# columns (H:N) on several rows across the ALC, ARM, BSM, CRP, CUL, GSM and
# LTW sheets with freshly recalculated values. Some cells that previously
# had no applicable value are populated for the first time, and some cells
# that are no longer applicable are cleared.

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 123.07692
$ws.Range("I5").Value = 121
$ws.Range("J5").Value = 125.5
$ws.Range("K5").Value = 121
$ws.Range("L5").Value = 125.5
$ws.Range("M5").Value = -6
$ws.Range("N5").Value = -355.5
# Row 6
$ws.Range("H6").Value = 578.875
$ws.Range("I6").Value = 377.5
$ws.Range("K6").Value = 1132.5
$ws.Range("M6").Value = -1020.5
# Row 9
$ws.Range("H9").Value = 213.72728
$ws.Range("I9").Value = 181.5
$ws.Range("K9").Value = 181.5
$ws.Range("M9").Value = -12.5
# Row 17
$ws.Range("H17").Value = 1860.3182
$ws.Range("J17").Value = 1860.3182
$ws.Range("L17").Value = 5580.9546
$ws.Range("N17").Value = -5916.9546
# Row 32
$ws.Range("H32").Value = 798.25
$ws.Range("I32").Value = 796
$ws.Range("J32").Value = 799
$ws.Range("K32").Value = 796
$ws.Range("L32").Value = 799
$ws.Range("M32").Value = -470
$ws.Range("N32").Value = -1451
# Row 40
$ws.Range("H40").Value = 6204.4287
$ws.Range("I40").Value = 5366.8
$ws.Range("K40").Value = 5366.8
$ws.Range("M40").Value = -5191.8
# Row 53
$ws.Range("H53").Value = 321
$ws.Range("I53").Value = 368.2
$ws.Range("J53").Value = 85
$ws.Range("K53").Value = 368.2
$ws.Range("L53").Value = 85
$ws.Range("M53").Value = 268.8
$ws.Range("N53").Value = -1359
# Row 69
$ws.Range("H69").Value = 7220.409
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 7220.409
$ws.Range("K69").Value = 0
$ws.Range("L69").ClearContents()
$ws.Range("M69").Value = 21661.227
$ws.Range("N69").Value = -23409.227
# Row 72
$ws.Range("H72").Value = 7220.409
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 7220.409
$ws.Range("K72").Value = 0
$ws.Range("L72").ClearContents()
$ws.Range("M72").Value = 64983.681
$ws.Range("N72").Value = -73719.681
# Row 92
$ws.Range("H92").Value = 107.125
$ws.Range("I92").Value = 120.333336
$ws.Range("J92").Value = 67.5
$ws.Range("K92").Value = 120.333336
$ws.Range("L92").Value = 67.5
$ws.Range("M92").Value = 1127.666664
$ws.Range("N92").Value = -2563.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 272.91306
$ws.Range("I4").Value = 154.6842
$ws.Range("J4").Value = 834.5
$ws.Range("K4").Value = 154.6842
$ws.Range("L4").Value = 834.5
$ws.Range("M4").Value = -38.6842
$ws.Range("N4").Value = -1066.5
# Row 61
$ws.Range("H61").Value = 5998.75
$ws.Range("I61").Value = 5998.75
$ws.Range("K61").Value = 5998.75
$ws.Range("M61").Value = -5786.75
# Row 136
$ws.Range("H136").Value = 5998.75
$ws.Range("I136").Value = 5998.75
$ws.Range("K136").Value = 17996.25
$ws.Range("M136").Value = -15446.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3241.2727
$ws.Range("J20").Value = 3465.4
$ws.Range("L20").Value = 3465.4
$ws.Range("N20").Value = -3959.4
# Row 94
$ws.Range("H94").Value = 749
$ws.Range("I94").Value = 748.5
$ws.Range("K94").Value = 748.5
$ws.Range("M94").Value = -297.5
# Row 99
$ws.Range("H99").Value = 1432.6
$ws.Range("I99").Value = 1037.8462
$ws.Range("K99").Value = 1037.8462
$ws.Range("M99").Value = 460.1538
# Row 107
$ws.Range("H107").Value = 4885.8335
$ws.Range("I107").Value = 2829.4167
$ws.Range("J107").Value = 8998.666999999999
$ws.Range("K107").Value = 2829.4167
$ws.Range("L107").Value = 8998.666999999999
$ws.Range("M107").Value = -909.4167000000002
$ws.Range("N107").Value = -12838.667

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 17500
$ws.Range("I4").Value = 17500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 17500
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -17388
# Row 31
$ws.Range("H31").Value = 4564.8535
$ws.Range("I31").Value = 2440.3794
$ws.Range("J31").Value = 9699
$ws.Range("K31").Value = 2440.3794
$ws.Range("L31").Value = 9699
$ws.Range("M31").Value = -2145.3794
$ws.Range("N31").Value = -10289
# Row 34
$ws.Range("H34").Value = 4564.8535
$ws.Range("I34").Value = 2440.3794
$ws.Range("J34").Value = 9699
$ws.Range("K34").Value = 2440.3794
$ws.Range("L34").Value = 9699
$ws.Range("M34").Value = -2238.3794
$ws.Range("N34").Value = -10103
# Row 58
$ws.Range("H58").Value = 5631.8335
$ws.Range("I58").Value = 4949
$ws.Range("K58").Value = 4949
$ws.Range("M58").Value = -4746
# Row 69
$ws.Range("H69").Value = 11750
$ws.Range("I69").Value = 4000
$ws.Range("J69").Value = 35000
$ws.Range("K69").Value = 4000
$ws.Range("L69").Value = 35000
$ws.Range("M69").Value = -3251
$ws.Range("N69").Value = -36498
# Row 72
$ws.Range("H72").Value = 11750
$ws.Range("I72").Value = 4000
$ws.Range("J72").Value = 35000
$ws.Range("K72").Value = 12000
$ws.Range("L72").Value = 105000
$ws.Range("M72").Value = -8256
$ws.Range("N72").Value = -112488
# Row 136
$ws.Range("H136").Value = 5631.8335
$ws.Range("I136").Value = 4949
$ws.Range("K136").Value = 14847
$ws.Range("M136").Value = -12297

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 211438.3
$ws.Range("J4").Value = 18730.334
$ws.Range("L4").Value = 56191.00199999999
$ws.Range("N4").Value = -56415.00199999999
# Row 59
$ws.Range("H59").Value = 933.3333
$ws.Range("I59").Value = 933.3333
$ws.Range("K59").Value = 2799.9999
$ws.Range("M59").Value = -2259.9999
# Row 61
$ws.Range("H61").Value = 200
$ws.Range("I61").Value = 200
$ws.Range("K61").Value = 600
$ws.Range("M61").Value = -385
# Row 107
$ws.Range("H107").Value = 1999.4
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1999.4
$ws.Range("K107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("M107").Value = 5998.200000000001
$ws.Range("N107").Value = -9838.200000000001
# Row 115
$ws.Range("H115").Value = 2499.5
$ws.Range("I115").Value = 2000
$ws.Range("K115").Value = 6000
$ws.Range("M115").Value = -4825

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
# Row 73
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
# Row 80
$ws.Range("H80").Value = 3396.6667
$ws.Range("I80").Value = 3345
$ws.Range("K80").Value = 3345
$ws.Range("M80").Value = -2347
# Row 83
$ws.Range("H83").Value = 3396.6667
$ws.Range("I83").Value = 3345
$ws.Range("K83").Value = 16725
$ws.Range("M83").Value = -11733
# Row 97
$ws.Range("H97").Value = 946
$ws.Range("I97").Value = 758.6667
$ws.Range("J97").Value = 1133.3334
$ws.Range("K97").Value = 758.6667
$ws.Range("L97").Value = 1133.3334
$ws.Range("M97").Value = -262.6667
$ws.Range("N97").Value = -2125.3334

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 1564.3
$ws.Range("I55").Value = 1486.625
$ws.Range("K55").Value = 1486.625
$ws.Range("M55").Value = -1313.625
